$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: holiday excel path config entry key (description filled in later, see below)
$ws.Range("A8").Value = "holidayExcelPath"

# Rows 9-16: SCADA tag / virtual weather station name pairs for the states
$ws.Range("A9").Value = "WRLDCMP.SCADA1.A0047000"
$ws.Range("B9").Value = "WR virtual weather station name"

$ws.Range("A10").Value = "WRLDCMP.SCADA1.A0046980"
$ws.Range("B10").Value = "mah virtual weather station name"

$ws.Range("A11").Value = "WRLDCMP.SCADA1.A0046957"
$ws.Range("B11").Value = "Gujarat  virtual weather station name"

$ws.Range("A12").Value = "WRLDCMP.SCADA1.A0046978"
$ws.Range("B12").Value = "MP  virtual weather station name"

$ws.Range("A13").Value = "WRLDCMP.SCADA1.A0046945"
$ws.Range("B13").Value = "Chatt  virtual weather station name"

$ws.Range("A14").Value = "WRLDCMP.SCADA1.A0046962"
$ws.Range("B14").Value = "Goa  virtual weather station name"

$ws.Range("A15").Value = "WRLDCMP.SCADA1.A0046948"
$ws.Range("B15").Value = "DD  virtual weather station name"

$ws.Range("A16").Value = "WRLDCMP.SCADA1.A0046953"
$ws.Range("B16").Value = "DNH  virtual weather station name"

# Description for the holiday excel path key, added last (so its shared string is appended at the end)
$ws.Range("B8").Value = "holiday excel folder path"
# Match the same "broken hyperlink" formatting used by the other description cells in column B
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selected cell to match the post-edit state (selection moved to B22 in the source)
$ws.Range("B22").Select()
